# Auto-generated edit script: applies updated FFXIV leve-profit market
# data values to the Alexander_Profits workbook, per scheduled runner diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 187.16667
$ws.Range("I2").Value = 187.16667
$ws.Range("K2").Value = 187.16667
$ws.Range("M2").Value = -74.16667000000001
$ws.Range("H19").Value = 376.75
$ws.Range("I19").Value = 413.4
$ws.Range("J19").Value = 350.57144
$ws.Range("K19").Value = 413.4
$ws.Range("L19").Value = 350.57144
$ws.Range("M19").Value = -238.4
$ws.Range("N19").Value = -700.5714399999999
$ws.Range("H82").Value = 15251.5
$ws.Range("I82").Value = 1003
$ws.Range("K82").Value = 3009
$ws.Range("M82").Value = -2603
$ws.Range("H85").Value = 15251.5
$ws.Range("I85").Value = 1003
$ws.Range("K85").Value = 3009
$ws.Range("M85").Value = -1605
$ws.Range("H112").Value = 2649.2856
$ws.Range("J112").Value = 3149.1304
$ws.Range("L112").Value = 9447.3912
$ws.Range("N112").Value = -11663.3912
$ws.Range("H129").Value = 639.4
$ws.Range("I129").Value = 450.90625
$ws.Range("J129").Value = 2650
$ws.Range("K129").Value = 1352.71875
$ws.Range("L129").Value = 7950
$ws.Range("M129").Value = 3647.28125
$ws.Range("N129").Value = -17950
$ws.Range("H137").Value = 3410393
$ws.Range("I137").Value = 1924384.5
$ws.Range("J137").Value = 5556850
$ws.Range("K137").Value = 5773153.5
$ws.Range("L137").Value = 16670550
$ws.Range("M137").Value = -5770603.5
$ws.Range("N137").Value = -16675650
$ws.Range("H140").Value = 55015
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 55015
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 55015
$ws.Range("N140").Value = -65375
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 2349170.8
$ws.Range("I141").Value = 1281.9524
$ws.Range("J141").Value = 5749561.5
$ws.Range("K141").Value = 3845.857199999999
$ws.Range("L141").Value = 17248684.5
$ws.Range("M141").Value = 1334.142800000001
$ws.Range("N141").Value = -17259044.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 33721.43
$ws.Range("J7").Value = 33721.43
$ws.Range("L7").Value = 33721.43
$ws.Range("N7").Value = -33949.43
$ws.Range("H32").Value = 2677912.5
$ws.Range("I32").Value = 4620.8047
$ws.Range("J32").Value = 20945406
$ws.Range("K32").Value = 4620.8047
$ws.Range("L32").Value = 20945406
$ws.Range("M32").Value = -4333.8047
$ws.Range("N32").Value = -20945980
$ws.Range("H122").Value = 2392.6667
$ws.Range("I122").Value = 2271.2
$ws.Range("K122").Value = 6813.599999999999
$ws.Range("M122").Value = -4363.599999999999
$ws.Range("H124").Value = 21621.45
$ws.Range("J124").Value = 21621.45
$ws.Range("L124").Value = 21621.45
$ws.Range("N124").Value = -31441.45
$ws.Range("H132").Value = 84124.42999999999
$ws.Range("I132").Value = 99419.96000000001
$ws.Range("J132").Value = 6117.2
$ws.Range("K132").Value = 298259.88
$ws.Range("L132").Value = 18351.6
$ws.Range("M132").Value = -295729.88
$ws.Range("N132").Value = -23411.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 40092
$ws.Range("J122").Value = 40092
$ws.Range("L122").Value = 40092
$ws.Range("N122").Value = -49892
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1682.8649
$ws.Range("I31").Value = 1302.5294
$ws.Range("J31").Value = 5993.3335
$ws.Range("K31").Value = 1302.5294
$ws.Range("L31").Value = 5993.3335
$ws.Range("M31").Value = -1007.5294
$ws.Range("N31").Value = -6583.3335
$ws.Range("H34").Value = 1682.8649
$ws.Range("I34").Value = 1302.5294
$ws.Range("J34").Value = 5993.3335
$ws.Range("K34").Value = 1302.5294
$ws.Range("L34").Value = 5993.3335
$ws.Range("M34").Value = -1100.5294
$ws.Range("N34").Value = -6397.3335
$ws.Range("H99").Value = 1799.75
$ws.Range("I99").Value = 1799.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1799.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -301.75
$ws.Range("N99").ClearContents()
$ws.Range("H120").Value = 49986.668
$ws.Range("J120").Value = 49986.668
$ws.Range("L120").Value = 49986.668
$ws.Range("N120").Value = -57244.668
$ws.Range("H122").Value = 2113.6155
$ws.Range("I122").Value = 1171.4286
$ws.Range("J122").Value = 3212.8333
$ws.Range("K122").Value = 3514.2858
$ws.Range("L122").Value = 9638.499899999999
$ws.Range("M122").Value = -1064.2858
$ws.Range("N122").Value = -14538.4999
$ws.Range("H126").Value = 1799.75
$ws.Range("I126").Value = 1799.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5399.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2929.25
$ws.Range("N126").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 46777.164
$ws.Range("I122").Value = 352.29413
$ws.Range("J122").Value = 57442.336
$ws.Range("K122").Value = 3170.64717
$ws.Range("L122").Value = 516981.024
$ws.Range("M122").Value = -720.6471700000002
$ws.Range("N122").Value = -521881.024
$ws.Range("H137").Value = 2005.35
$ws.Range("I137").Value = 1757.8572
$ws.Range("J137").Value = 2582.8333
$ws.Range("K137").Value = 5273.571599999999
$ws.Range("L137").Value = 7748.499899999999
$ws.Range("M137").Value = -173.5715999999993
$ws.Range("N137").Value = -17948.4999
$ws.Range("H141").Value = 10669.3125
$ws.Range("I141").Value = 4412.1113
$ws.Range("K141").Value = 13236.3339
$ws.Range("M141").Value = -8056.333899999998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5980
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 5980
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 5980
$ws.Range("N102").Value = -9224
$ws.Range("M102").ClearContents()
$ws.Range("H113").Value = 1266.6428
$ws.Range("I113").Value = 1115.3529
$ws.Range("K113").Value = 1115.3529
$ws.Range("M113").Value = 1054.6471
$ws.Range("H126").Value = 6433.5293
$ws.Range("I126").Value = 3570
$ws.Range("J126").Value = 8438
$ws.Range("K126").Value = 10710
$ws.Range("L126").Value = 25314
$ws.Range("M126").Value = -8240
$ws.Range("N126").Value = -30254
$ws.Range("H141").Value = 64940.6
$ws.Range("J141").Value = 64940.6
$ws.Range("L141").Value = 64940.6
$ws.Range("N141").Value = -75300.60000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14287786
$ws.Range("I7").Value = 20001500
$ws.Range("K7").Value = 20001500
$ws.Range("M7").Value = -20001388
$ws.Range("H46").Value = 983.175
$ws.Range("I46").Value = 1046.1052
$ws.Range("J46").Value = 926.2381
$ws.Range("K46").Value = 1046.1052
$ws.Range("L46").Value = 926.2381
$ws.Range("M46").Value = -858.1052
$ws.Range("N46").Value = -1302.2381
$ws.Range("H118").Value = 31729.2
$ws.Range("J118").Value = 31729.2
$ws.Range("L118").Value = 31729.2
$ws.Range("N118").Value = -35043.2
$ws.Range("H126").Value = 14287786
$ws.Range("I126").Value = 20001500
$ws.Range("K126").Value = 60004500
$ws.Range("M126").Value = -60002030
$ws.Range("H127").Value = 55087
$ws.Range("J127").Value = 55087
$ws.Range("L127").Value = 55087
$ws.Range("N127").Value = -65007
$ws.Range("H132").Value = 1788.7188
$ws.Range("I132").Value = 1699.8937
$ws.Range("J132").Value = 2034.2941
$ws.Range("K132").Value = 5099.6811
$ws.Range("L132").Value = 6102.8823
$ws.Range("M132").Value = -2569.6811
$ws.Range("N132").Value = -11162.8823
$ws.Range("H136").Value = 1873.1818
$ws.Range("I136").Value = 1846.7307
$ws.Range("J136").Value = 1971.4286
$ws.Range("K136").Value = 5540.1921
$ws.Range("L136").Value = 5914.2858
$ws.Range("M136").Value = -2990.1921
$ws.Range("N136").Value = -11014.2858

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 18994
$ws.Range("J42").Value = 18994
$ws.Range("L42").Value = 18994
$ws.Range("N42").Value = -19750
$ws.Range("H126").Value = 1695.1666
$ws.Range("I126").Value = 1230.25
$ws.Range("J126").Value = 2625
$ws.Range("K126").Value = 3690.75
$ws.Range("L126").Value = 7875
$ws.Range("M126").Value = -1220.75
$ws.Range("N126").Value = -12815

Write-Host "Applied all Sheets updates via scheduled runner."
